$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold prices as plain text (e.g. "19.29", "1.50",
# "0.0504"). Assigning such a numeric-looking string straight to
# .Value lets Excel coerce it to a Number (dropping significant
# trailing zeros, e.g. "1.50" -> 1.5). Mark just the D cells we are
# about to rewrite as Text first so the literal string is preserved,
# without touching the NumberFormat of any other cell.
$priceCells = @("D2","D3","D5","D8","D10","D12","D13","D15","D16","D17","D18","D19","D20","D22","D23","D24","D25","D27","D29","D33","D34","D35","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.653.68'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '1.620.49'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '214.48'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("E6").Value = '  +0.57%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.246'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("D10").Value = '19.29'
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '1.850.46'
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("D13").Value = '1.618.16'
$ws.Range("E13").Value = '  +2.35%  '
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").Value = '64.70'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = '0.511'
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("D17").Value = '26.686.65'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").Value = '233.13'
$ws.Range("E18").Value = '  +9.44%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0729'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '7.70'
$ws.Range("E20").Value = '  +4.58%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '4.39'
$ws.Range("E22").Value = '  +3.11%  '
$ws.Range("D23").Value = '2.26'
$ws.Range("E23").Value = '  +4.34%  '
$ws.Range("D24").Value = '9.08'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("D25").Value = '145.74'
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '7.05'
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("E28").Value = '  +2.20%  '
$ws.Range("D29").Value = '15.61'
$ws.Range("E29").Value = '  +2.82%  '
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("E31").Value = '  +1.11%  '
$ws.Range("E32").Value = '  +1.71%  '
$ws.Range("D33").Value = '1.462.06'
$ws.Range("E33").Value = '  +9.23%  '
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  +1.95%  '
$ws.Range("D35").Value = '2.43'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("E36").Value = '  +1.68%  '
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = '0.837'
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("D40").Value = '5.92'
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '2.21'
$ws.Range("E42").Value = '  +3.17%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '0.952'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '1.760.55'
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("D45").Value = '0.764'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").Value = '61.86'
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").Value = '88.45'
$ws.Range("E47").Value = '  +3.22%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.50'
$ws.Range("E49").Value = '  +1.58%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0504'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.0963'
$ws.Range("E51").Value = '  -1.61%  '
